$wb = $excel.ActiveWorkbook

# Column F holds "想去人数" (number of people interested). The generated
# data was refreshed, bumping several counts. Column index F = 6.

$updates = @{
    "展览"     = @{ 3=1016; 4=247; 5=13; 7=734; 8=254; 11=408; 12=216; 13=83; 14=853; 16=1998; 17=492; 18=7502; 19=552; 20=521; 21=60; 22=95; 24=225 }
    "演出"     = @{ 8=123; 10=6 }
    "本地生活" = @{ 2=5530 }
    "全部类型" = @{ 3=5530; 7=1016; 10=247; 11=13; 13=734; 14=254; 18=408; 19=216; 21=83; 23=853; 25=123; 26=1998; 27=492; 28=7502; 30=6; 31=552; 32=521; 33=60; 34=95; 37=225 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 6).Value = $rows[$r]
    }
}
